$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-04-21"

# Update the header label for the 2022 column (shared string)
$ws.Range("I1").Value = "2022 (through 04-21)"

# Update April 2022 value (row 5) and recalc-able Total (row 14)
$ws.Range("I5").Value = 91
$ws.Range("I14").Value = 527
